$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 900
$ws.Range("I12").Value = 1300
$ws.Range("J12").Value = 100
$ws.Range("K12").Value = 1300
$ws.Range("L12").Value = 100
$ws.Range("M12").Value = -1130
$ws.Range("N12").Value = -440
$ws.Range("H39").Value = 346.25
$ws.Range("I39").Value = 281.25
$ws.Range("J39").Value = 411.25
$ws.Range("K39").Value = 843.75
$ws.Range("L39").Value = 1233.75
$ws.Range("M39").Value = -547.75
$ws.Range("N39").Value = -1825.75
$ws.Range("H100").Value = 13264361
$ws.Range("I100").Value = 22865042
$ws.Range("J100").Value = 63425
$ws.Range("K100").Value = 22865042
$ws.Range("L100").Value = 63425
$ws.Range("M100").Value = -22864501
$ws.Range("N100").Value = -64507
$ws.Range("H112").Value = 2289.0588
$ws.Range("I112").Value = 1674
$ws.Range("J112").Value = 2478.3076
$ws.Range("K112").Value = 5022
$ws.Range("L112").Value = 7434.9228
$ws.Range("M112").Value = -3914
$ws.Range("N112").Value = -9650.9228
$ws.Range("H134").Value = 89995.53999999999
$ws.Range("J134").Value = 89995.53999999999
$ws.Range("L134").Value = 89995.53999999999
$ws.Range("N134").Value = -100135.54
$ws.Range("H138").Value = 214721.95
$ws.Range("I138").Value = 399839.28
$ws.Range("J138").Value = 4071.2068
$ws.Range("K138").Value = 1199517.84
$ws.Range("L138").Value = 12213.6204
$ws.Range("M138").Value = -1194377.84
$ws.Range("N138").Value = -22493.6204
$ws.Range("H139").Value = 79633.28999999999
$ws.Range("J139").Value = 87905.5
$ws.Range("L139").Value = 87905.5
$ws.Range("N139").Value = -98185.5
$ws.Range("H140").Value = 85460
$ws.Range("J140").Value = 85460
$ws.Range("L140").Value = 85460
$ws.Range("N140").Value = -95820

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H36").Value = 32500
$ws.Range("I36").Value = 32500
$ws.Range("K36").Value = 32500
$ws.Range("M36").Value = -32154
$ws.Range("H61").Value = 5824.409
$ws.Range("I61").Value = 6094.6284
$ws.Range("J61").Value = 4773.5557
$ws.Range("K61").Value = 6094.6284
$ws.Range("L61").Value = 4773.5557
$ws.Range("M61").Value = -5882.6284
$ws.Range("N61").Value = -5197.5557
$ws.Range("H110").Value = 1948.08
$ws.Range("I110").Value = 1570.4762
$ws.Range("K110").Value = 1570.4762
$ws.Range("M110").Value = 474.5237999999999
$ws.Range("H132").Value = 2177.804
$ws.Range("I132").Value = 1717.6389
$ws.Range("J132").Value = 3282.2
$ws.Range("K132").Value = 5152.9167
$ws.Range("L132").Value = 9846.599999999999
$ws.Range("M132").Value = -2622.9167
$ws.Range("N132").Value = -14906.6
$ws.Range("H133").Value = 72630
$ws.Range("J133").Value = 72630
$ws.Range("L133").Value = 72630
$ws.Range("N133").Value = -77690
$ws.Range("H136").Value = 5824.409
$ws.Range("I136").Value = 6094.6284
$ws.Range("J136").Value = 4773.5557
$ws.Range("K136").Value = 18283.8852
$ws.Range("L136").Value = 14320.6671
$ws.Range("M136").Value = -15733.8852
$ws.Range("N136").Value = -19420.6671
$ws.Range("H140").Value = 70000
$ws.Range("J140").Value = 70000
$ws.Range("L140").Value = 70000
$ws.Range("N140").Value = -80360
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = $null

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 11330
$ws.Range("I8").Value = 4000
$ws.Range("J8").Value = 14995
$ws.Range("K8").Value = 4000
$ws.Range("L8").Value = 14995
$ws.Range("M8").Value = -3860
$ws.Range("N8").Value = -15275
$ws.Range("H94").Value = 2898.08
$ws.Range("J94").Value = 5583.857
$ws.Range("L94").Value = 5583.857
$ws.Range("N94").Value = -6485.857
$ws.Range("H134").Value = 5691.0347
$ws.Range("I134").Value = 6958.409
$ws.Range("K134").Value = 20875.227
$ws.Range("M134").Value = -18340.227
$ws.Range("H140").Value = 78057.25
$ws.Range("J140").Value = 78057.25
$ws.Range("L140").Value = 78057.25
$ws.Range("N140").Value = -88417.25
$ws.Range("H141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("L141").Value = 0
$ws.Range("N141").Value = $null

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 429.25
$ws.Range("I19").Value = 429.25
$ws.Range("K19").Value = 429.25
$ws.Range("M19").Value = -259.25
$ws.Range("H24").Value = 429.25
$ws.Range("I24").Value = 429.25
$ws.Range("K24").Value = 429.25
$ws.Range("M24").Value = -259.25
$ws.Range("H31").Value = 6836.9653
$ws.Range("I31").Value = 7012.4165
$ws.Range("J31").Value = 5994.8
$ws.Range("K31").Value = 7012.4165
$ws.Range("L31").Value = 5994.8
$ws.Range("M31").Value = -6717.4165
$ws.Range("N31").Value = -6584.8
$ws.Range("H34").Value = 6836.9653
$ws.Range("I34").Value = 7012.4165
$ws.Range("J34").Value = 5994.8
$ws.Range("K34").Value = 7012.4165
$ws.Range("L34").Value = 5994.8
$ws.Range("M34").Value = -6810.4165
$ws.Range("N34").Value = -6398.8
$ws.Range("H47").Value = 46999
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").Value = $null
$ws.Range("H122").Value = 10911.846
$ws.Range("I122").Value = 18540.572
$ws.Range("K122").Value = 55621.716
$ws.Range("M122").Value = -53171.716
$ws.Range("H132").Value = 1592.4286
$ws.Range("I132").Value = 1452.1578
$ws.Range("J132").Value = 2925
$ws.Range("K132").Value = 4356.4734
$ws.Range("L132").Value = 8775
$ws.Range("M132").Value = -1826.4734
$ws.Range("N132").Value = -13835
$ws.Range("H133").Value = 80000.5
$ws.Range("J133").Value = 80000.5
$ws.Range("L133").Value = 80000.5
$ws.Range("N133").Value = -85060.5
$ws.Range("H138").Value = 63000
$ws.Range("J138").Value = 63000
$ws.Range("L138").Value = 63000
$ws.Range("N138").Value = -73280
$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("M139").Value = $null
$ws.Range("N139").Value = $null
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").Value = $null
$ws.Range("H141").Value = 312586.1
$ws.Range("I141").Value = 90000
$ws.Range("K141").Value = 90000
$ws.Range("M141").Value = -84820

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 675.2222
$ws.Range("J32").Value = 753.375
$ws.Range("L32").Value = 2260.125
$ws.Range("N32").Value = -2826.125
$ws.Range("H98").Value = 1442.2222
$ws.Range("J98").Value = 1074.3334
$ws.Range("L98").Value = 3223.0002
$ws.Range("N98").Value = -6219.0002

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 11465.363
$ws.Range("I122").Value = 9570.333000000001
$ws.Range("J122").Value = 13739.4
$ws.Range("K122").Value = 28710.999
$ws.Range("L122").Value = 41218.2
$ws.Range("M122").Value = -26260.999
$ws.Range("N122").Value = -46118.2
$ws.Range("H132").Value = 2259.6428
$ws.Range("I132").Value = 2216.262
$ws.Range("J132").Value = 2389.7856
$ws.Range("K132").Value = 6648.786
$ws.Range("L132").Value = 7169.3568
$ws.Range("M132").Value = -4118.786
$ws.Range("N132").Value = -12229.3568
$ws.Range("H140").Value = 91354.42999999999
$ws.Range("J140").Value = 91354.42999999999
$ws.Range("L140").Value = 91354.42999999999
$ws.Range("N140").Value = -101714.43

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6268.8945
$ws.Range("I22").Value = 11577.444
$ws.Range("K22").Value = 11577.444
$ws.Range("M22").Value = -11282.444
$ws.Range("H27").Value = 6268.8945
$ws.Range("I27").Value = 11577.444
$ws.Range("K27").Value = 11577.444
$ws.Range("M27").Value = -11470.444
$ws.Range("H50").Value = 62436.5
$ws.Range("I50").Value = 54877
$ws.Range("K50").Value = 54877
$ws.Range("M50").Value = -54240
$ws.Range("H55").Value = 1806.8334
$ws.Range("I55").Value = 409.6
$ws.Range("J55").Value = 2804.8572
$ws.Range("K55").Value = 409.6
$ws.Range("L55").Value = 2804.8572
$ws.Range("M55").Value = -236.6
$ws.Range("N55").Value = -3150.8572
$ws.Range("H100").Value = 4457.7646
$ws.Range("I100").Value = 3026
$ws.Range("K100").Value = 3026
$ws.Range("M100").Value = -2485
$ws.Range("H122").Value = 5617.4443
$ws.Range("I122").Value = 6115.2856
$ws.Range("J122").Value = 3875
$ws.Range("K122").Value = 18345.8568
$ws.Range("L122").Value = 11625
$ws.Range("M122").Value = -15895.8568
$ws.Range("N122").Value = -16525
$ws.Range("H132").Value = 554231.5600000001
$ws.Range("I132").Value = 785600.2
$ws.Range("K132").Value = 2356800.6
$ws.Range("M132").Value = -2354270.6

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H75").Value = 26952.5
$ws.Range("J75").Value = 32655
$ws.Range("L75").Value = 32655
$ws.Range("N75").Value = -34527
$ws.Range("H78").Value = 26952.5
$ws.Range("J78").Value = 32655
$ws.Range("L78").Value = 97965
$ws.Range("N78").Value = -107325
$ws.Range("H107").Value = 17009.895
$ws.Range("I107").Value = 1288.2222
$ws.Range("K107").Value = 3864.6666
$ws.Range("M107").Value = -1944.6666
$ws.Range("H122").Value = 4121.7637
$ws.Range("I122").Value = 1709.4839
$ws.Range("J122").Value = 7237.625
$ws.Range("K122").Value = 5128.4517
$ws.Range("L122").Value = 21712.875
$ws.Range("M122").Value = -2678.4517
$ws.Range("N122").Value = -26612.875
$ws.Range("H132").Value = 9843.788
$ws.Range("I132").Value = 12058.243
$ws.Range("J132").Value = 4381.467
$ws.Range("K132").Value = 36174.729
$ws.Range("L132").Value = 13144.401
$ws.Range("M132").Value = -33644.729
$ws.Range("N132").Value = -18204.401
$ws.Range("H137").Value = 80140.86
$ws.Range("J137").Value = 80140.86
$ws.Range("L137").Value = 80140.86
$ws.Range("N137").Value = -90340.86

Write-Host "Applied all changes"